$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column R header date (Excel serial date 45394 = 2024-04-12)
$ws.Range("R2").Value = 45394

# New text cells in column R (set in the order the shared strings must be
# appended: Bericht versenden(22), Finalisierung(23), Anordnung...(24))
$ws.Range("R3").Value = "Bericht versenden"
$ws.Range("R5").Value = "Finalisierung "
$ws.Range("R4").Value = "Anordnung von Maschinen endern"
